# Update the "想去人数" (number of interested people) column F values
# on the "展览" (Exhibition) and "全部类型" (All Types) sheets to reflect
# the latest scraped counts.

$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> row number -> new F value
$updates = @{
    "展览" = @{
        2  = 5529
        4  = 3
        5  = 359
        7  = 2
        9  = 30
        10 = 59
        12 = 139
        13 = 339
        14 = 434
        15 = 3043
        16 = 3
        18 = 1659
    }
    "全部类型" = @{
        2  = 5529
        4  = 3
        5  = 359
        7  = 2
        10 = 30
        11 = 59
        13 = 139
        14 = 339
        15 = 434
        16 = 3043
        17 = 3
        19 = 1659
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($r in $rows.Keys) {
        $ws.Cells.Item($r, 6).Value = $rows[$r]
    }
}
